$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.5
$ws.Range("K2").Value = 2.05
$ws.Range("AG2").Value = 12
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 101
$ws.Range("AQ2").Value = 26
$ws.Range("AW2").Value = 8.5
$ws.Range("AZ2").Value = 251

# Row 3 updates
$ws.Range("G3").Value = 4.2
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 2.15
$ws.Range("O3").Value = 1.73
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 3.5
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 1.75
$ws.Range("T3").Value = 2.05
$ws.Range("AC3").Value = 4.75
$ws.Range("AF3").Value = 126
$ws.Range("AJ3").Value = 21
$ws.Range("AK3").Value = 26
$ws.Range("AX3").Value = 15
$ws.Range("BB3").Value = 451
